# Apply the data update described in the diff.
# The rows 298-306 (Poland Ekstraklasa fixtures) had their match-data
# (everything from column B "id" through column AD "PL_AhUnder", i.e.
# skipping the row-counter in column A) reassigned between rows, following
# two independent cyclic rotations:
#   298 -> 299 -> 300 -> 298                 (3-row cycle)
#   301 -> 302 -> 303 -> 304 -> 305 -> 306 -> 301   (6-row cycle)
# Concretely, the NEW content of each row equals the OLD content of the
# "source" row below it (wrapping around within its cycle).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values to copy FROM, as they existed
# before this edit).
$mapping = @{
    298 = 299
    299 = 300
    300 = 298
    301 = 302
    302 = 303
    303 = 304
    304 = 305
    305 = 306
    306 = 301
}

# Snapshot the original values for columns B:AD (id .. PL_AhUnder) for every
# row involved, before any writes happen, so the permutation is applied
# consistently regardless of write order.
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $addr = "B" + $row + ":AD" + $row
    $rng = $ws.Range($addr)
    $snapshot[$row] = $rng.Value2
}

# Now write each destination row using the snapshot of its source row.
foreach ($row in $mapping.Keys) {
    $srcRow = $mapping[$row]
    $addr = "B" + $row + ":AD" + $row
    $destRng = $ws.Range($addr)
    $destRng.Value2 = $snapshot[$srcRow]
}
